{"js": "// CV edits (paragraph-scoped searches so repeated substrings like\n// \"coding \" never cross-match another paragraph):\n//  1. Contact line: \"sam | 0701563048 | samuel.mburu@tracom.co.ke\"\n//     -> \"sam mburu | 0700000000 | sam@gmail\"\n//  2. \"i love coding\" -> \"i love coding in python, java\"\n//  3. Bold job-title run \"tracom services  \" -> \"tracom \"\n//  4. Trailing run \"amazing content \" -> \"programming \"\n//  5. First bullet \"coding \" -> \"coding\" (drop trailing space)\n//  6. Second bullet \"programmer\" -> \"coding  after coding\"\n\n// Single-run paragraphs: clear() + insertText() replaces the run cleanly\n// (no leftover run-level artifacts from the old text). Safe here because\n// clear() empties the WHOLE paragraph, and each of these paragraphs\n// contains exactly one run anyway.\nasync function replaceWholeParagraph(context, paragraph, expectedOldText, newText) {\n  paragraph.load(\"text\");\n  await context.sync();\n  if (paragraph.text !== expectedOldText) {\n    throw new Error(\n      \"Paragraph text mismatch: expected \" + JSON.stringify(expectedOldText) +\n        \", found \" + JSON.stringify(paragraph.text)\n    );\n  }\n  paragraph.clear();\n  await context.sync();\n  paragraph.insertText(newText, \"Start\");\n  await context.sync();\n}\n\n// Multi-run paragraphs: only touch the matched sub-range so sibling runs\n// (different bold/italic formatting, breaks, etc.) are left completely\n// untouched. (paragraph.clear()/range.clear() operate at whole-paragraph\n// granularity, so they are NOT used here.)\nasync function replaceInParagraph(context, paragraph, findText, newText) {\n  const results = paragraph.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(findText) +\n        \" in paragraph, found \" + results.items.length\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Paragraph 1: contact line (single run).\nawait replaceWholeParagraph(\n  context,\n  items[1],\n  \"sam | 0701563048 | samuel.mburu@tracom.co.ke\",\n  \"sam mburu | 0700000000 | sam@gmail\"\n);\n\n// Paragraph 3: \"About me\" body line (single run).\nawait replaceWholeParagraph(context, items[3], \"i love coding\", \"i love coding in python, java\");\n\n// Paragraph 5: \"Work Experience\" entry (bold title run + italic date run\n// with a line break + plain trailing run) -- edit sub-ranges only.\nawait replaceInParagraph(context, items[5], \"tracom services  \", \"tracom \");\nawait replaceInParagraph(context, items[5], \"amazing content \", \"programming \");\n\n// Paragraph 7: first Skills bullet (single run).\nawait replaceWholeParagraph(context, items[7], \"coding \", \"coding\");\n\n// Paragraph 8: second Skills bullet (single run).\nawait replaceWholeParagraph(context, items[8], \"programmer\", \"coding  after coding\");\n", "ps1": "# CV edits via Word COM interop. Each Find/Replace is scoped to the\n# specific paragraph's Range so that repeated substrings (e.g. \"coding \")\n# never cross-match a different paragraph, and only the intended run's\n# text changes (sibling runs such as the bold title / italic date in the\n# \"Work Experience\" entry are left completely untouched).\n\n$d = $word.ActiveDocument\n\nfunction Replace-InRange($range, $findText, $replaceText) {\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $findText\n    $range.Find.Replacement.Text = $replaceText\n    $result = $range.Find.Execute(\n        [ref]$findText,   # FindText\n        [ref]$true,       # MatchCase\n        [ref]$false,      # MatchWholeWord\n        [ref]$false,      # MatchWildcards\n        [ref]$false,      # MatchSoundsLike\n        [ref]$false,      # MatchAllWordForms\n        [ref]$true,       # Forward\n        [ref]1,           # Wrap (wdFindContinue)\n        [ref]$null,       # Format\n        [ref]$replaceText,# ReplaceWith\n        [ref]2            # Replace (wdReplaceAll)\n    )\n    if (-not $result) {\n        throw \"Find/Replace failed for '$findText'\"\n    }\n}\n\n# Paragraph 2: contact line.\nReplace-InRange $d.Paragraphs(2).Range \"sam | 0701563048 | samuel.mburu@tracom.co.ke\" \"sam mburu | 0700000000 | sam@gmail\"\n\n# Paragraph 4: \"About me\" body line.\nReplace-InRange $d.Paragraphs(4).Range \"i love coding\" \"i love coding in python, java\"\n\n# Paragraph 6: \"Work Experience\" entry (bold title run + trailing run).\nReplace-InRange $d.Paragraphs(6).Range \"tracom services  \" \"tracom \"\nReplace-InRange $d.Paragraphs(6).Range \"amazing content \" \"programming \"\n\n# Paragraph 8: first Skills bullet.\nReplace-InRange $d.Paragraphs(8).Range \"coding \" \"coding\"\n\n# Paragraph 9: second Skills bullet.\nReplace-InRange $d.Paragraphs(9).Range \"programmer\" \"coding  after coding\"\n"}
